$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every data cell as inline/shared text, even for values
# that look numeric or like a percentage ("315.22", "3.57%", ...), and none of the
# touched cells carry an explicit cell style (default "General" / style index 0).
# A plain `.Value =` assignment lets Excel auto-convert such numeric-looking
# strings into real numbers/percentages (and stamps a NumberFormat on the cell),
# so instead we briefly force Text format to keep the literal string, then put the
# cell style back to the workbook default so the XML stays style-index-0, exactly
# like the untouched cells around it.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '315.22'
Set-TextValue "E2" '3.57%'
Set-TextValue "D3" '35.44'
Set-TextValue "E3" '-0.38%'
Set-TextValue "D4" '5.103'
Set-TextValue "E4" '1.00%'
Set-TextValue "D5" '0.08169'
Set-TextValue "E5" '3.72%'
Set-TextValue "D6" '2.061'
Set-TextValue "E6" '-3.39%'
Set-TextValue "B7" 'GateToken'
Set-TextValue "C7" 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D7" '4.141'
Set-TextValue "E7" '0.04%'
Set-TextValue "B8" 'KuCoinToken'
Set-TextValue "C8" 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue "D8" '7.942'
Set-TextValue "E8" '0.02%'
Set-TextValue "B9" 'MXToken'
Set-TextValue "C9" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D9" '0.9316'
Set-TextValue "E9" '0.76%'
Set-TextValue "B10" 'LiechtensteinCryptoassetsExchange'
Set-TextValue "C10" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D10" '0.1038'
Set-TextValue "E10" '6.45%'
Set-TextValue "B11" 'WazirX'
Set-TextValue "C11" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D11" '0.1921'
Set-TextValue "E11" '4.58%'
Set-TextValue "B12" 'MandalaExchangeToken'
Set-TextValue "C12" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D12" '0.09073'
Set-TextValue "E12" '5.43%'
Set-TextValue "B13" 'BitrueCoin'
Set-TextValue "C13" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D13" '0.03599'
Set-TextValue "E13" '0.45%'
Set-TextValue "B14" 'BitMartToken'
Set-TextValue "C14" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D14" '0.09889'
Set-TextValue "E14" '-0.34%'
Set-TextValue "B15" 'BitForexToken'
Set-TextValue "C15" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D15" '0.001437'
Set-TextValue "E15" '0.37%'
Set-TextValue "B16" 'TigerCash'
Set-TextValue "C16" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D16" '0.005680'
Set-TextValue "E16" '-0.65%'
Set-TextValue "B17" 'LEO'
Set-TextValue "C17" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D17" '3.468'
Set-TextValue "E17" '-0.15%'
Set-TextValue "E18" '3.65%'
Set-TextValue "D19" '0.3407'
Set-TextValue "E19" '0.98%'
Set-TextValue "E20" '-4.06%'
Set-TextValue "D21" '5.095'
Set-TextValue "E21" '-1.30%'
Set-TextValue "D23" '0.04549'
Set-TextValue "E23" '-0.14%'
Set-TextValue "E24" '0.78%'
Set-TextValue "D25" '0.004796'
Set-TextValue "E25" '-0.78%'
Set-TextValue "D26" '0.0001250'
Set-TextValue "E26" '-3.88%'
Set-TextValue "D27" '0.0004500'
Set-TextValue "E27" '-5.26%'
Set-TextValue "D39" '0.01980'
Set-TextValue "E39" '6.89%'
Set-TextValue "D40" '0.04967'
Set-TextValue "E40" '5.06%'
Set-TextValue "D41" '0.007604'
Set-TextValue "E41" '-2.35%'
Set-TextValue "D42" '0.1382'
Set-TextValue "E42" '-0.27%'
Set-TextValue "D43" '0.007866'
Set-TextValue "E43" '1.43%'
Set-TextValue "D44" '0.002260'
Set-TextValue "E44" '4.46%'
Set-TextValue "D45" '0.01177'
Set-TextValue "E45" '3.90%'
Set-TextValue "D46" '0.00006624'
Set-TextValue "E46" '5.20%'
Set-TextValue "D47" '0.00000000750'
Set-TextValue "E47" '-0.02%'
Set-TextValue "D48" '64.07'
Set-TextValue "E48" '26.35%'
Set-TextValue "D49" '0.001700'
Set-TextValue "E49" '-10.54%'
Set-TextValue "D50" '0.00002100'
Set-TextValue "E50" '-0.02%'
Set-TextValue "D51" '0.0002000'
Set-TextValue "E51" '-0.02%'
